$wb = $excel.ActiveWorkbook

# Add a new row to the "prepare" sheet with a "file" entry pointing to the
# newly added source workbook.
$prepare = $wb.Worksheets.Item("prepare")
$prepare.Cells.Item(7, 1).Value = 6
$prepare.Cells.Item(7, 2).Value = "file"
$prepare.Cells.Item(7, 3).Value = "source/source_03.xlsx"
$prepare.Columns.Item(3).ColumnWidth = 22.166666666666668

# Make the "source" sheet the active/selected sheet instead of "prepare".
$source = $wb.Worksheets.Item("source")
$source.Activate()
